$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "usuario"
$ws.Range("B6").Value = 1234

$ws.Range("B6").Select()
